$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 2 (A2:T2) gets new cluster-name values and recomputed metrics.
# Rows 3-7 are brand-new rows appended below, following the same column layout.

$data = @(
    @("ECs",  "Hp", "Itgam", "ECs", 1, 0.3333333333333333, 33.102689,          99.30806699999999, 0.6442155376576643, 0.6442155376576643, 1, 0.3333333333333333, 66.38494866666666, 199.154846,  0.4171200956172241, 0.4171200956172242, 2197.520309993631, 19777.68278994268, 0.2687152466658664, 0.2687152466658664),
    @("ECs",  "Hp", "Itgam", "M2",  1, 0.3333333333333333, 33.102689,          99.30806699999999, 0.6442155376576643, 0.6442155376576643, 3, 1,                  92.76573566666667, 278.297207,  0.5828799043827758, 0.5828799043827758, 3070.795297629874, 27637.15767866887, 0.3755002909917979, 0.3755002909917979),
    @("FAPs", "Hp", "Itgam", "ECs", 2, 0.6666666666666666, 1.719635666666666,  5.158906999999999, 0.03346604306305637, 0.03346604306305637, 1, 0.3333333333333333, 66.38494866666666, 199.154846,  0.4171200956172241, 0.4171200956172242, 114.1579254570357, 1027.421329113322, 0.01395935908239222, 0.01395935908239222),
    @("FAPs", "Hp", "Itgam", "M2",  2, 0.6666666666666666, 1.719635666666666,  5.158906999999999, 0.03346604306305637, 0.03346604306305637, 3, 1,                  92.76573566666667, 278.297207,  0.5828799043827758, 0.5828799043827758, 159.5232676969721, 1435.709409272749, 0.01950668398066416, 0.01950668398066416),
    @("M2",   "Hp", "Itgam", "ECs", 3, 1,                  16.56216866666666,  49.68650599999999, 0.3223184192792793, 0.3223184192792793, 1, 0.3333333333333333, 66.38494866666666, 199.154846,  0.4171200956172241, 0.4171200956172242, 1099.478716745341, 9895.308450708075, 0.1344454898689655, 0.1344454898689655),
    @("M2",   "Hp", "Itgam", "M2",  3, 1,                  16.56216866666666,  49.68650599999999, 0.3223184192792793, 0.3223184192792793, 3, 1,                  92.76573566666667, 278.297207,  0.5828799043827758, 0.5828799043827758, 1536.401760598749, 13827.61584538874, 0.1878729294103137, 0.1878729294103137)
)

$row = 2
foreach ($record in $data) {
    $col = 1
    foreach ($value in $record) {
        $ws.Cells.Item($row, $col).Value = $value
        $col++
    }
    $row++
}
